$wb = $excel.ActiveWorkbook

# --- Fill in the Bump_3G simulation results (rows 6-15, columns D, E, F) ---
$ws = $wb.Worksheets.Item("Bump_3G")

$ws.Range("D6").Value = -0.256
$ws.Range("E6").Value = -0.00917
$ws.Range("F6").Value = 15

$ws.Range("D7").Value = -0.158
$ws.Range("E7").Value = 0.0135
$ws.Range("F7").Value = 15

$ws.Range("D8").Value = -0.158
$ws.Range("E8").Formula = "=2.72*10^-6"
$ws.Range("F8").Value = 12.3

$ws.Range("D9").Value = -0.157
$ws.Range("E9").Value = -0.00526
$ws.Range("F9").Value = 12.3

$ws.Range("D10").Value = -0.161
$ws.Range("E10").Value = -0.00741
$ws.Range("F10").Value = 8.91

$ws.Range("D11").Value = -0.163
$ws.Range("E11").Value = 0.00174
$ws.Range("F11").Value = 8.91

$ws.Range("D12").Value = 0.321
$ws.Range("E12").Value = -0.0146
$ws.Range("F12").Value = 3.32

$ws.Range("D13").Value = 0.331
$ws.Range("E13").Value = 0.0164
$ws.Range("F13").Value = 3.34

$ws.Range("D14").Value = 0.0652
$ws.Range("E14").Value = 0.0069
$ws.Range("F14").Value = 1.61

$ws.Range("D15").Value = 0.0904
$ws.Range("E15").Value = -0.00281
$ws.Range("F15").Value = 1.61

# --- Move the active tab / selection from "MAX SPEED" to "Bump_3G" ---
$wsMaxSpeed = $wb.Worksheets.Item("MAX SPEED")
$wsMaxSpeed.Activate()
$wsMaxSpeed.Range("G24").Select()

$ws.Activate()
$ws.Range("F16").Select()
